$wb = $excel.ActiveWorkbook

# --- "runs" sheet: update run-range + sim time ---
$runs = $wb.Worksheets.Item("runs")
$runs.Range("B1").Value = 3
$runs.Range("B2").Value = 3
$runs.Range("B3").Value = 200

# --- "params" sheet: update D6 and move its stored selection there ---
$params = $wb.Worksheets.Item("params")
$params.Range("D6").Value = 1
$params.Range("D6").Select()

# Restore "runs" as the active/visible tab (it was active before this edit)
$runs.Activate()
